$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet already contains 4 identical 27-row "clinical" schedule blocks
# for groups B1A (rows 2-28), B1B (29-55), B1C (56-82) and B1D (83-109).
# A 5th block for group B1E (rows 110-136) needs to be appended, using the
# exact same dates/times/durations as the B1D block.

# 1) Copy the values (text/numbers) of the B1D block into the new B1E rows.
$valueSource = $ws.Range("A83:G109")
$valueSource.Copy()
$destination = $ws.Range("A110:G136")
$destination.PasteSpecial(-4163)   # xlPasteValues - keeps "1".."27" as text, 180 as number

# 2) Copy the cell formatting from the B1A block (rows 2-28), which starts on
#    an even row just like the destination (row 110), so the alternating
#    row-banding styles line up correctly with the new rows.
$formatSource = $ws.Range("A2:G28")
$formatSource.Copy()
$destination.PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = 0

# 3) Update the Group column for the new rows to "B1E".
$ws.Range("B110:B136").Value = "B1E"
